# Auto-derived from the upstream OOXML diff.
$wb = $excel.ActiveWorkbook

# --- 1) Metadata sheet: refresh URL + Date -------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/audit-retention"
$wsMeta.Range("B8").Value = "2025-08-20T10:40:04+01:00"

# --- 2) Elements sheet: re-fit column widths ------------------------------
$ws = $wb.Worksheets.Item("Elements")
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666  # target width 16.41796875
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666  # target width 16.41796875
$ws.Columns.Item(3).Hidden = $true
$ws.Columns.Item(3).ColumnWidth = 9.0  # target width 9.79296875
$ws.Columns.Item(4).Hidden = $true
$ws.Columns.Item(4).ColumnWidth = 6.166666666666667  # target width 7.046875
$ws.Columns.Item(5).ColumnWidth = 4.5  # target width 5.30078125
$ws.Columns.Item(6).ColumnWidth = 3.1666666666666665  # target width 3.953125
$ws.Columns.Item(7).ColumnWidth = 3.5  # target width 4.265625
$ws.Columns.Item(8).ColumnWidth = 11.833333333333334  # target width 12.6875
$ws.Columns.Item(9).ColumnWidth = 9.666666666666666  # target width 10.51171875
$ws.Columns.Item(10).ColumnWidth = 19.833333333333332  # target width 20.703125
$ws.Columns.Item(11).ColumnWidth = 7.5  # target width 8.3984375
$ws.Columns.Item(12).ColumnWidth = 99.83333333333333  # target width 100.703125
$ws.Columns.Item(13).ColumnWidth = 99.83333333333333  # target width 100.703125
$ws.Columns.Item(14).ColumnWidth = 99.83333333333333  # target width 100.703125
$ws.Columns.Item(15).ColumnWidth = 11.5  # target width 12.26171875
$ws.Columns.Item(16).ColumnWidth = 19.833333333333332  # target width 20.703125
$ws.Columns.Item(17).ColumnWidth = 19.833333333333332  # target width 20.703125
$ws.Columns.Item(18).ColumnWidth = 19.833333333333332  # target width 20.703125
$ws.Columns.Item(19).ColumnWidth = 19.833333333333332  # target width 20.703125
$ws.Columns.Item(20).ColumnWidth = 7.0  # target width 7.80078125
$ws.Columns.Item(21).ColumnWidth = 12.833333333333334  # target width 13.609375
$ws.Columns.Item(22).ColumnWidth = 13.166666666666666  # target width 13.91796875
$ws.Columns.Item(23).ColumnWidth = 14.166666666666666  # target width 15.01171875
$ws.Columns.Item(24).ColumnWidth = 13.833333333333334  # target width 14.62890625
$ws.Columns.Item(25).ColumnWidth = 16.166666666666668  # target width 17.08203125
$ws.Columns.Item(26).ColumnWidth = 14.333333333333334  # target width 15.18359375
$ws.Columns.Item(27).ColumnWidth = 4.166666666666667  # target width 5.07421875
$ws.Columns.Item(28).ColumnWidth = 17.166666666666668  # target width 17.98046875
$ws.Columns.Item(29).ColumnWidth = 33.666666666666664  # target width 34.578125
$ws.Columns.Item(30).ColumnWidth = 12.666666666666666  # target width 13.54296875
$ws.Columns.Item(31).Hidden = $true
$ws.Columns.Item(31).ColumnWidth = 10.5  # target width 11.3203125
$ws.Columns.Item(32).Hidden = $true
$ws.Columns.Item(32).ColumnWidth = 14.166666666666666  # target width 15.046875
$ws.Columns.Item(33).Hidden = $true
$ws.Columns.Item(33).ColumnWidth = 7.333333333333333  # target width 8.22265625
$ws.Columns.Item(34).ColumnWidth = 7.666666666666667  # target width 8.53125
$ws.Columns.Item(35).ColumnWidth = 99.83333333333333  # target width 100.703125
$ws.Columns.Item(37).ColumnWidth = 18.666666666666668  # target width 19.5625

Write-Host "audit-retention metadata + column widths updated"
